$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, shifting existing rows 164..279 down to 165..280
$ws.Rows("164").Insert()

# Populate the newly inserted row 164 with the new weekly data point
$ws.Range("A164").Value = 4
$ws.Range("B164").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C164").Value = "Los Lagos"
$ws.Range("D164").Value = 44574
$ws.Range("E164").Value = 10
$ws.Range("F164").Value = 100114013
$ws.Range("G164").Value = "Zanahoria"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 250
$ws.Range("K164").Value = 12000
$ws.Range("L164").Value = 14000
$ws.Range("M164").Value = 13200
$ws.Range("N164").Value = "`$/saco 20 kilos"
$ws.Range("O164").Value = "Región de Ñuble"
$ws.Range("P164").Value = 660
$ws.Range("Q164").Value = 20
$ws.Range("R164").Value = "Hortaliza"
